# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.548.56'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.815.43'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.53%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.57%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.002'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '306.10'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.73%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4528'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.72%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3603'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.61%  '
$ws.Range("E9").Value = '  +2.95%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07101'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.34%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8986'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07785'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.40'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.833.94'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.285'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.327'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '85.39'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.78%  '
$ws.Range("E18").Value = '  -0.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008561'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.597.29'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.25'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.978'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("E24").Value = '  +1.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.012.81'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.16%  '
$ws.Range("E26").Value = '  -1.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '150.69'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.83'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.50%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.063'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '112.58'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.865'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.63%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08710'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.56%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.124'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.87%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7523'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.754'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.455'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.112'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.073'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01935'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05116'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.905'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.77%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5106'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.762'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1509'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.58%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.070'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.58%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4741'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.10%  '
$ws.Range("E47").Value = '  -0.73%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.07'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.62%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '101.03'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.53%  '
$ws.Range("E50").Value = '  -0.46%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05984'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.03%  '
